$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Admin User should be able to Add Subject Categories", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Position Catagories Add Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "FAILED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "FAILED", "chrome"),
    @("Human Resources  Position Catagories Add Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  Attestations Functionality", "PASSED", "chrome"),
    @("Human Resources  PositionCatagories Delete Functionality", "PASSED", "chrome")
)

$r = 10
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
